# (Closes #25029) Removed extra bracket from cheatsheet code example. (#25032)
#
# The "Create DataFrame with a MultiIndex" code sample contains a stray
# trailing closing parenthesis:
#     names=['n','v'])))
# It should read:
#     names=['n','v']))
# Find the shape/run containing the typo (wherever it lives) and fix the
# text in place, leaving every other run/paragraph untouched.

$p = $ppt.ActivePresentation

$needle = "names=['n','v'])))"
$target = $null

foreach ($slide in $p.Slides) {
    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text.IndexOf($needle) -ge 0) {
                $target = $shape
                break
            }
        }
    }
    if ($target -ne $null) {
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the shape containing the cheatsheet code sample."
}

$tr = $target.TextFrame.TextRange
$found = $tr.Find("'])))")
$found.Text = "']))"
